$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.816.02'
$ws.Range("E2").Value = '  +4.87%  '

$ws.Range("D3").Value = '1.617.72'
$ws.Range("E3").Value = '  +4.45%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.995'
$ws.Range("E4").Value = '  -0.45%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.90'
$ws.Range("E5").Value = '  +2.05%  '

$ws.Range("E6").Value = '  +7.36%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.995'
$ws.Range("E7").Value = '  -0.50%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '27.04'
$ws.Range("E8").Value = '  +12.77%  '

$ws.Range("E9").Value = '  +4.00%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0600'
$ws.Range("E10").Value = '  +3.13%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0915'
$ws.Range("E11").Value = '  +2.93%  '

$ws.Range("D12").Value = '1.849.24'
$ws.Range("E12").Value = '  +4.40%  '

$ws.Range("D13").Value = '1.606.04'
$ws.Range("E13").Value = '  +3.74%  '

$ws.Range("D14").Value = '29.850.35'
$ws.Range("E14").Value = '  +5.07%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.542'
$ws.Range("E15").Value = '  +6.66%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.78'
$ws.Range("E16").Value = '  +4.33%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '246.38'
$ws.Range("E17").Value = '  +7.61%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.65'
$ws.Range("E19").Value = '  +4.72%  '

$ws.Range("D20").Value = '0.0₃0698'
$ws.Range("E20").Value = '  +4.05%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.995'
$ws.Range("E21").Value = '  -0.48%  '

$ws.Range("E22").Value = '  +5.09%  '

$ws.Range("E23").Value = '  +4.57%  '

$ws.Range("E24").Value = '  +5.30%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '156.17'
$ws.Range("E25").Value = '  +3.27%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.44'
$ws.Range("E26").Value = '  +4.63%  '

$ws.Range("E27").Value = '  +6.16%  '

$ws.Range("E28").Value = '  +3.70%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.995'
$ws.Range("E29").Value = '  -0.54%  '

$ws.Range("E30").Value = '  +1.85%  '

$ws.Range("E31").Value = '  +0.70%  '

$ws.Range("E32").Value = '  +3.48%  '

$ws.Range("D33").Value = '1.448.97'
$ws.Range("E33").Value = '  +4.74%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.13'
$ws.Range("E34").Value = '  +4.61%  '

$ws.Range("E35").Value = '  -0.80%  '

$ws.Range("E36").Value = '  +4.08%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.86'
$ws.Range("E37").Value = '  +11.24%  '

$ws.Range("E38").Value = '  +0.21%  '

$ws.Range("E39").Value = '  +3.51%  '

$ws.Range("B40").Value = 'BitcoinSV'
$ws.Range("C40").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '57.10'
$ws.Range("E40").Value = '  +31.88%  '

$ws.Range("B41").Value = 'ImmutableX'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.542'
$ws.Range("E41").Value = '  +6.39%  '

$ws.Range("E42").Value = '  +2.47%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '69.62'
$ws.Range("E43").Value = '  +12.66%  '

$ws.Range("E44").Value = '  +4.48%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.994'
$ws.Range("E45").Value = '  -0.53%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0468'
$ws.Range("E46").Value = '  +2.32%  '

$ws.Range("E47").Value = '  +0.73%  '

$ws.Range("D48").Value = '1.760.38'
$ws.Range("E48").Value = '  +4.60%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '86.74'
$ws.Range("E49").Value = '  +1.35%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.837'
$ws.Range("E50").Value = '  -4.51%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0521'
$ws.Range("E51").Value = '  +2.01%  '
